$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "Sheet11"

$new.Range("A1").Value = "Power (dbm)"
$new.Range("C1").Value = "End time"
$new.Range("B1").Value = "Start Time"

$new.Range("A2").Value = -7
$new.Range("A3").Value = -6
$new.Range("A4").Value = -5

$new.Range("B2:C8").NumberFormat = "h:mm:ss"

$new.Range("B2").Value = 0.58472222222222225
$new.Range("C2").Value = 0.58576388888888886
$new.Range("B3").Value = 0.58611111111111114
$new.Range("C3").Value = 0.58750000000000002
$new.Range("B4").Value = 0.58819444444444446
$new.Range("C4").Value = 0.58969907407407407

$new.Range("A7").Select()

Write-Host "done"
